$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("breakfast_theroost")

# Column D = LocalIngredients: clear the leftover "Egg"/"Egg, Sausage" placeholder
# values on rows 2-6 and replace with "NA" to mark them explicitly empty.
$ws.Range("D2:D6").Value = "NA"

# Column F = nutritionLabel: replace the unfinished "waffles" values with a
# generic "placeholder" across all data rows.
$ws.Range("F2:F9").Value = "placeholder"

$ws.Range("G23").Select()
